$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the now-unused "Calculations" worksheet entirely.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Calculations").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 2. "About" sheet: replace the old "Notes:" block with a short "Source:"
#    note referencing CPL's consultation with the American Forest Foundation.
# ---------------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Drop everything from row 4 down (the old multi-paragraph "Notes:" text).
$about.Rows("4:10").Delete() | Out-Null

# A3 becomes "Source:" (was "Notes:"), with the source note itself in B3.
$about.Cells.Item(3, 1).Value = "Source:"
$about.Cells.Item(3, 2).Value = "consultation with American Forest Foundation"

# Recreate the small amount of leftover formatting that remains below it.
$about.Cells.Item(5, 1).HorizontalAlignment = -4131   # xlLeft
$about.Rows(6).RowHeight = 14.45

# The old "Notes:" paragraphs were the only thing using the hyperlink-style
# run/cell style in this workbook; it is no longer referenced anywhere, so
# drop the now-unused named cell style along with it.
$wb.Styles.Item("Hiperlink").Delete()

# ---------------------------------------------------------------------------
# 3. "CApULAbIFM" sheet: CO2 Abated (g) is now computed via a formula
#    instead of the previous hard-coded 0.
# ---------------------------------------------------------------------------
$data = $wb.Worksheets.Item("CApULAbIFM")
$data.Range("B2").Formula = "=1.5*10^6"
